$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

$ws.Rows.Item(95).Insert(-4121)
$tbl.Resize($ws.Range("A8:K143"))

$ws.Range("A94:K94").Copy()
$ws.Range("A95:K95").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G95").Formula = $ws.Range("G94").Formula()
$ws.Range("G143").Formula = $ws.Range("G142").Formula()

$ws.Range("C93").Value = 1.25
$ws.Range("B94").Value = "SL(1-0-0)"
$ws.Range("H94").Value = 1
$ws.Range("K94").Value = 45087

$ws.Range("B95").Value = "SL(1-0-0)"
$ws.Range("H95").Value = 1
$ws.Range("K95").Value = 45094

# Copy exact style (incl. date number format) from an existing date-styled K cell
$ws.Range("K15").Copy()
$ws.Range("K94").PasteSpecial(-4122)
$ws.Range("K95").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K94").Value = 45087
$ws.Range("K95").Value = 45094

$excel.CalculateFull()
